$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": the dc47a288 entry moves from
# "In Translation" to "Ready for handoff" with fresh handoff timestamps,
# on the Overview sheet and on each per-locale sheet (zh-cn, de-de).

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-13-17 14:13:48"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-17 14:13:44"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-17 14:13:48"
